# "[2 24] update ex14"
#
# The sheet currently holds an Excel Table ("TabelaOsoby") covering A1:C4:
#   row 1 -> header labels (Imie / Wiek / Miasto)
#   rows 2-4 -> the actual people data (Jan/29/Warszawa, Anna/24/Wroclaw, Tomek/35/Legnica)
#
# Target state: the table goes away (converted back to a plain range) and the
# header row is removed entirely, leaving only the data in rows 2-4 (dimension
# becomes A2:C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the Excel Table back into a normal range (drops the table
# definition/autofilter/style banding, keeps the cell data in place).
if ($ws.ListObjects.Count -gt 0) {
    $ws.ListObjects.Item(1).Unlist()
}

# Remove the header row's content so the sheet starts at row 2.
$ws.Range("A1:C1").ClearContents()
